$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: new data row -------------------------------------------------
$ws.Range("A5").Value = "DSAlgo085"
$ws.Range("A5").BorderAround(1)

$ws.Range("B5").Value = "Hello@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:Hello@gmail.com")

$ws.Range("C5").Value = "Hello098@"
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:Hello098@")

# --- selection ------------------------------------------------------------
$ws.Range("C7").Select()

Write-Host "done"
